$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 2 header dates (row 10) ---
$ws.Range("C10").Value = 43486
$ws.Range("E10").Value = 43487
$ws.Range("G10").Value = 43488
$ws.Range("I10").Value = 43489
$ws.Range("K10").Value = 43490

# --- "Планируемые часы работы" (intended/planned hours) block ---
# Morning session (row 11)
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0.46875
$ws.Range("F11").Value = 8
$ws.Range("G11").Value = 0.48958333333333331
$ws.Range("H11").Value = 8
$ws.Range("I11").Value = 0.42708333333333331
$ws.Range("J11").Value = 8

# Afternoon session (row 12)
$ws.Range("E12").Value = 0.80208333333333337
$ws.Range("G12").Value = 0.82291666666666663
$ws.Range("I12").Value = 0.76041666666666663

# --- "Фактические часы работы" (actual today hours) block ---
# Morning session (row 14)
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0.46875
$ws.Range("F14").Value = 8

# Afternoon session (row 15)
$ws.Range("E15").Value = 0.80208333333333337

# Update the view: select K11 as active cell (also clears the stale
# topLeftCell="A4" scroll position left over from the previous selection).
$ws.Range("K11").Select()
